# Prefix each worksheet's name onto the "Step..." / command names held in
# column A (data rows only - row 1 is the header row containing "Name").
#
# This applies to the per-protocol sheets added after the core journey
# sheets (price1, price2, discount1, discount2, free1, free2, nomoney1,
# nomoney2, noppv1, noppv2, card1, card2, nosex1, nosex2, offtopic1,
# offtopic2, real1, real2, voice1, voice2, customyes1, customyes2,
# customno1, customno2, done1, done2, cumcontrol, dickpic, boosters).

$wb = $excel.ActiveWorkbook

# Sheets that should NOT receive a prefix (the original journey sheets).
$excluded = @("ZackJourney", "MeetupRedirect", "NRWaves", "PersonalZack", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    $sheetName = $ws.Name
    if ($excluded -contains $sheetName) {
        continue
    }

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value2
        if ($null -ne $val -and $val -ne "") {
            $prefix = $sheetName + " "
            if (-not $val.ToString().StartsWith($prefix)) {
                $cell.Value = $prefix + $val
            }
        }
    }
}
